$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NPM validation data / names per new sorting by jurusan
$ws.Range("B5").Value = "H1A021099"
$ws.Range("A3").Value = "Muhammad ssed"
$ws.Range("B3").Value = "G1A021066"

# Update active selection
$ws.Range("C14").Select()
